# Facilitator guidelines - Ants Problem: Swahili (Kenya) -> English translations,
# plus the document's default proofing language locale tag: sw-KE -> sw-TZ.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # wdFindContinue = 1 (Wrap), wdReplaceAll = 2
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Kichwa cha Video" "Video Title"
Replace-Text "Mada" "Topic"
Replace-Text "Malengo" "Aim(s)"
Replace-Text "Urefu" "Length"
Replace-Text "Mahali pa Kambi" "Camp Location"
Replace-Text "Wawezeshaji" "Facilitators"
Replace-Text "N. ya wanafunzi" "N. of students"
Replace-Text "Tarehe" "Date"
Replace-Text "Rasilimali" "Resources"
Replace-Text "inahitajika" "needed"
Replace-Text "Maandalizi" "Preparations"
Replace-Text "Muda wa video" "Video time"
Replace-Text "Mwezeshaji anafanya nini" "What facilitator does"
Replace-Text "Wanachofanya wanafunzi" "What learners do"
Replace-Text "Utangulizi Mkuu wa Video ya VMC" "General VMC Video Introduction"
Replace-Text "Utangulizi wa Video" "Video Introduction"
Replace-Text "Kitendawili" "Riddle"
# Occurs twice in the document; wdReplaceAll handles both occurrences.
Replace-Text "Kusaidia mchakato, kuchochea mawazo" "Assist the process, provoke thoughts"
Replace-Text "Suluhisho" "Solution"

# Document default language: Swahili (Kenya) -> Swahili (Tanzania).
# All body text uses the "Normal" style, so set the default/base style's
# proofing language accordingly.
$normal = $d.Styles("Normal")
$normal.LanguageID = "sw-TZ"
